# Implements commit 'feat: add 2022-Q1 data':
#   - insert a new '2022-Q1' worksheet (fund-holdings detail) right
#     before the existing '总计' (summary) worksheet
#   - prepend a matching '2022-Q1' row to the '总计' summary table
$wb = $excel.ActiveWorkbook

# ---- helpers -------------------------------------------------------
# Excel auto-coerces a pure-numeric-looking string into a Number when
# assigned through .Value. The source data stores fund codes / ratios
# as TEXT, so force text entry, then paste-special just the formats
# from a known-good template cell to avoid leaving a stray numFmt='@'
# style behind.
function Set-TextValue {
    param($Cell, $Value, $FormatSource)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $FormatSource.Copy()
    $Cell.PasteSpecial(-4122)
}

# Copy only the formatting (xlPasteFormats) from a template cell.
function Copy-Style {
    param($Cell, $FormatSource)
    $FormatSource.Copy()
    $Cell.PasteSpecial(-4122)
}

# 1) Insert the new '2022-Q1' sheet right before '总计' -------------
$totalBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalBefore)
$q1.Name = "2022-Q1"

# IMPORTANT: this COM shim resolves a worksheet object by its *index*
# at access time rather than binding to a stable object identity, so
# $totalBefore now (mis)reports as the sheet that ended up at its old
# index (i.e. the freshly inserted '2022-Q1' sheet). Re-fetch '总计'
# by name now that the sheet order has settled, and use that handle
# for everything below.
$total = $wb.Worksheets.Item("总计")

# Style templates: '总计' already carries the bold/bordered header
# style and the plain data-cell style this new sheet reuses.
$headerStyleSrc = $total.Range("B1")
$aColStyleSrc = $total.Range("A2")
$plainStyleSrc = $total.Range("B2")

# --- header row ---
$c = $q1.Cells.Item(1, 2)
Set-TextValue $c '基金代码' $headerStyleSrc
$c = $q1.Cells.Item(1, 3)
Set-TextValue $c '基金名称' $headerStyleSrc
$c = $q1.Cells.Item(1, 4)
Set-TextValue $c '基金规模' $headerStyleSrc
$c = $q1.Cells.Item(1, 5)
Set-TextValue $c '股票总仓位' $headerStyleSrc
$c = $q1.Cells.Item(1, 6)
Set-TextValue $c '仓位占比' $headerStyleSrc
$c = $q1.Cells.Item(1, 7)
Set-TextValue $c '持有市值(亿元)' $headerStyleSrc
$c = $q1.Cells.Item(1, 8)
Set-TextValue $c '仓位排名' $headerStyleSrc

# --- data rows ---
# row 2
$c = $q1.Cells.Item(2, 1)
$c.Value = 0
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(2, 2)
Set-TextValue $c '010695' $plainStyleSrc
$c = $q1.Cells.Item(2, 3)
Set-TextValue $c '华夏磐益一年定期开放混合' $plainStyleSrc
$c = $q1.Cells.Item(2, 4)
Set-TextValue $c '18.02' $plainStyleSrc
$c = $q1.Cells.Item(2, 5)
Set-TextValue $c '82.41' $plainStyleSrc
$c = $q1.Cells.Item(2, 6)
Set-TextValue $c '2.53' $plainStyleSrc
$c = $q1.Cells.Item(2, 7)
Set-TextValue $c '0.4559' $plainStyleSrc
$c = $q1.Cells.Item(2, 8)
$c.Value = 9
# row 3
$c = $q1.Cells.Item(3, 1)
$c.Value = 1
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(3, 2)
Set-TextValue $c '010861' $plainStyleSrc
$c = $q1.Cells.Item(3, 3)
Set-TextValue $c '长信企业优选一年持有期灵活配置混合' $plainStyleSrc
$c = $q1.Cells.Item(3, 4)
Set-TextValue $c '9.39' $plainStyleSrc
$c = $q1.Cells.Item(3, 5)
Set-TextValue $c '80.21' $plainStyleSrc
$c = $q1.Cells.Item(3, 6)
Set-TextValue $c '3.24' $plainStyleSrc
$c = $q1.Cells.Item(3, 7)
Set-TextValue $c '0.3042' $plainStyleSrc
$c = $q1.Cells.Item(3, 8)
$c.Value = 5
# row 4
$c = $q1.Cells.Item(4, 1)
$c.Value = 2
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(4, 2)
Set-TextValue $c '005589' $plainStyleSrc
$c = $q1.Cells.Item(4, 3)
Set-TextValue $c '长信企业精选两年定期开放灵活配置混合' $plainStyleSrc
$c = $q1.Cells.Item(4, 4)
Set-TextValue $c '5.84' $plainStyleSrc
$c = $q1.Cells.Item(4, 5)
Set-TextValue $c '79.99' $plainStyleSrc
$c = $q1.Cells.Item(4, 6)
Set-TextValue $c '3.22' $plainStyleSrc
$c = $q1.Cells.Item(4, 7)
Set-TextValue $c '0.1880' $plainStyleSrc
$c = $q1.Cells.Item(4, 8)
$c.Value = 5
# row 5
$c = $q1.Cells.Item(5, 1)
$c.Value = 3
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(5, 2)
Set-TextValue $c '162204' $plainStyleSrc
$c = $q1.Cells.Item(5, 3)
Set-TextValue $c '泰达宏利行业精选混合' $plainStyleSrc
$c = $q1.Cells.Item(5, 4)
Set-TextValue $c '7.28' $plainStyleSrc
$c = $q1.Cells.Item(5, 5)
Set-TextValue $c '75.56' $plainStyleSrc
$c = $q1.Cells.Item(5, 6)
Set-TextValue $c '1.87' $plainStyleSrc
$c = $q1.Cells.Item(5, 7)
Set-TextValue $c '0.1361' $plainStyleSrc
$c = $q1.Cells.Item(5, 8)
$c.Value = 10
# row 6
$c = $q1.Cells.Item(6, 1)
$c.Value = 4
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(6, 2)
Set-TextValue $c '003501' $plainStyleSrc
$c = $q1.Cells.Item(6, 3)
Set-TextValue $c '泰达宏利睿智稳健灵活配置混合' $plainStyleSrc
$c = $q1.Cells.Item(6, 4)
Set-TextValue $c '3.57' $plainStyleSrc
$c = $q1.Cells.Item(6, 5)
Set-TextValue $c '73.79' $plainStyleSrc
$c = $q1.Cells.Item(6, 6)
Set-TextValue $c '2.14' $plainStyleSrc
$c = $q1.Cells.Item(6, 7)
Set-TextValue $c '0.0764' $plainStyleSrc
$c = $q1.Cells.Item(6, 8)
$c.Value = 7
# row 7
$c = $q1.Cells.Item(7, 1)
$c.Value = 5
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(7, 2)
Set-TextValue $c '002863' $plainStyleSrc
$c = $q1.Cells.Item(7, 3)
Set-TextValue $c '金信深圳成长灵活配置混合' $plainStyleSrc
$c = $q1.Cells.Item(7, 4)
Set-TextValue $c '0.44' $plainStyleSrc
$c = $q1.Cells.Item(7, 5)
Set-TextValue $c '94.54' $plainStyleSrc
$c = $q1.Cells.Item(7, 6)
Set-TextValue $c '5.14' $plainStyleSrc
$c = $q1.Cells.Item(7, 7)
Set-TextValue $c '0.0226' $plainStyleSrc
$c = $q1.Cells.Item(7, 8)
$c.Value = 10
# row 8
$c = $q1.Cells.Item(8, 1)
$c.Value = 6
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(8, 2)
Set-TextValue $c '009128' $plainStyleSrc
$c = $q1.Cells.Item(8, 3)
Set-TextValue $c '明亚价值长青混合A' $plainStyleSrc
$c = $q1.Cells.Item(8, 4)
Set-TextValue $c '0.38' $plainStyleSrc
$c = $q1.Cells.Item(8, 5)
Set-TextValue $c '49.48' $plainStyleSrc
$c = $q1.Cells.Item(8, 6)
Set-TextValue $c '4.78' $plainStyleSrc
$c = $q1.Cells.Item(8, 7)
Set-TextValue $c '0.0182' $plainStyleSrc
$c = $q1.Cells.Item(8, 8)
$c.Value = 1
# row 9
$c = $q1.Cells.Item(9, 1)
$c.Value = 7
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(9, 2)
Set-TextValue $c '009658' $plainStyleSrc
$c = $q1.Cells.Item(9, 3)
Set-TextValue $c '汇丰晋信中小盘低波动策略股票A' $plainStyleSrc
$c = $q1.Cells.Item(9, 4)
Set-TextValue $c '0.98' $plainStyleSrc
$c = $q1.Cells.Item(9, 5)
Set-TextValue $c '86.56' $plainStyleSrc
$c = $q1.Cells.Item(9, 6)
Set-TextValue $c '1.09' $plainStyleSrc
$c = $q1.Cells.Item(9, 7)
Set-TextValue $c '0.0107' $plainStyleSrc
$c = $q1.Cells.Item(9, 8)
$c.Value = 8
# row 10
$c = $q1.Cells.Item(10, 1)
$c.Value = 8
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(10, 2)
Set-TextValue $c '008300' $plainStyleSrc
$c = $q1.Cells.Item(10, 3)
Set-TextValue $c '人保量化锐进混合A' $plainStyleSrc
$c = $q1.Cells.Item(10, 4)
Set-TextValue $c '0.20' $plainStyleSrc
$c = $q1.Cells.Item(10, 5)
Set-TextValue $c '92.95' $plainStyleSrc
$c = $q1.Cells.Item(10, 6)
Set-TextValue $c '2.31' $plainStyleSrc
$c = $q1.Cells.Item(10, 7)
Set-TextValue $c '0.0046' $plainStyleSrc
$c = $q1.Cells.Item(10, 8)
$c.Value = 7
# row 11
$c = $q1.Cells.Item(11, 1)
$c.Value = 9
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(11, 2)
Set-TextValue $c '007808' $plainStyleSrc
$c = $q1.Cells.Item(11, 3)
Set-TextValue $c '北信瑞丰量化优选灵活配置混合' $plainStyleSrc
$c = $q1.Cells.Item(11, 4)
Set-TextValue $c '0.24' $plainStyleSrc
$c = $q1.Cells.Item(11, 5)
Set-TextValue $c '79.84' $plainStyleSrc
$c = $q1.Cells.Item(11, 6)
Set-TextValue $c '1.06' $plainStyleSrc
$c = $q1.Cells.Item(11, 7)
Set-TextValue $c '0.0025' $plainStyleSrc
$c = $q1.Cells.Item(11, 8)
$c.Value = 8
# row 12
$c = $q1.Cells.Item(12, 1)
$c.Value = 10
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(12, 2)
Set-TextValue $c '008301' $plainStyleSrc
$c = $q1.Cells.Item(12, 3)
Set-TextValue $c '人保量化锐进混合C' $plainStyleSrc
$c = $q1.Cells.Item(12, 4)
Set-TextValue $c '0.06' $plainStyleSrc
$c = $q1.Cells.Item(12, 5)
Set-TextValue $c '92.95' $plainStyleSrc
$c = $q1.Cells.Item(12, 6)
Set-TextValue $c '2.31' $plainStyleSrc
$c = $q1.Cells.Item(12, 7)
Set-TextValue $c '0.0014' $plainStyleSrc
$c = $q1.Cells.Item(12, 8)
$c.Value = 7
# row 13
$c = $q1.Cells.Item(13, 1)
$c.Value = 11
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(13, 2)
Set-TextValue $c '009775' $plainStyleSrc
$c = $q1.Cells.Item(13, 3)
Set-TextValue $c '汇丰晋信中小盘低波动策略股票C' $plainStyleSrc
$c = $q1.Cells.Item(13, 4)
Set-TextValue $c '0.04' $plainStyleSrc
$c = $q1.Cells.Item(13, 5)
Set-TextValue $c '86.56' $plainStyleSrc
$c = $q1.Cells.Item(13, 6)
Set-TextValue $c '1.09' $plainStyleSrc
$c = $q1.Cells.Item(13, 7)
Set-TextValue $c '0.0004' $plainStyleSrc
$c = $q1.Cells.Item(13, 8)
$c.Value = 8
# row 14
$c = $q1.Cells.Item(14, 1)
$c.Value = 12
Copy-Style $c $aColStyleSrc
$c = $q1.Cells.Item(14, 2)
Set-TextValue $c '009129' $plainStyleSrc
$c = $q1.Cells.Item(14, 3)
Set-TextValue $c '明亚价值长青混合C' $plainStyleSrc
$c = $q1.Cells.Item(14, 4)
Set-TextValue $c '0.00' $plainStyleSrc
$c = $q1.Cells.Item(14, 5)
Set-TextValue $c '49.48' $plainStyleSrc
$c = $q1.Cells.Item(14, 6)
Set-TextValue $c '4.78' $plainStyleSrc
$c = $q1.Cells.Item(14, 7)
$c.Value = 0
$c = $q1.Cells.Item(14, 8)
$c.Value = 1

# 2) Prepend the 2022-Q1 row to '总计', shifting existing rows down -
# NOTE: this shim's `.Value` GETTER is unreliable (returns a stub
# descriptor instead of the real scalar); `.Value2` reads correctly,
# so use that when copying existing cell contents downward.
for ($r = 6; $r -ge 2; $r--) {
    $srcA = $total.Cells.Item($r, 1)
    $srcB = $total.Cells.Item($r, 2)
    $srcC = $total.Cells.Item($r, 3)
    $srcD = $total.Cells.Item($r, 4)
    $dstA = $total.Cells.Item($r + 1, 1)
    $dstB = $total.Cells.Item($r + 1, 2)
    $dstC = $total.Cells.Item($r + 1, 3)
    $dstD = $total.Cells.Item($r + 1, 4)
    # column A is the running index, so it increments by one as the
    # row moves down (row that used to be index N-1 becomes index N)
    $dstA.Value = $srcA.Value2 + 1
    Copy-Style $dstA $srcA
    $dstB.Value = $srcB.Value2
    $dstC.Value = $srcC.Value2
    $dstD.Value = $srcD.Value2
}

# New top row: index 0, '2022-Q1', 13 holdings, 1.22 亿元
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 13
$total.Cells.Item(2,4).Value = 1.22
